$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("prevalence2018")
$ws3.Range("C20:C36").Value = 0.049
$ws3.Range("C37:C46").Value = 0.062
$ws3.Range("C47:C56").Value = 0.056
$ws3.Range("C57:C82").Value = 0.068

$ws4 = $wb.Worksheets.Item("incidence2018_plus")
$ws4.Range("C20:C82").Formula = "=prevalence2018!C20/10"

$ws6 = $wb.Worksheets.Item("data")
$ws6.Range("D5:F8").Select()

$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("F20:F24").Select()

$ws4.Activate()
$ws4.Range("E20").Select()
